$wb = $excel.ActiveWorkbook

# Metadata sheet tracks Sample / Mineralogy for each analysis. A new
# analysis (SRM-88b, a dolomite) was run, so record its metadata in the
# next empty row.
$ws = $wb.Worksheets.Item("Metadata")
[void]$ws.Activate()

$ws.Cells.Item(2, 1).Value = "SRM-88b"
$ws.Cells.Item(2, 2).Value = "Dolomite"

[void]$ws.Range("B2").Select()
